$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.293105840682983
$ws.Range("B1").Value = 2.399503707885742
$ws.Range("C1").Value = 2.748632669448853
$ws.Range("D1").Value = 4.356668949127197
$ws.Range("E1").Value = 4.871416568756104
